$d = $word.ActiveDocument
$d.TrackRevisions = $false

# ---------------------------------------------------------------------------
# Change 1: paragraph "Número racional, recta numérica. " (keywords paragraph)
#   - re-split the runs ("n" / "úmero " / "racional,recta" / " " / "numérica")
#   - wrap "racional,recta" in a proofErr spellStart/spellEnd pair
#   - drop the trailing period
#   - move the "_GoBack" bookmark to the end of this paragraph
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute(
    "Número racional, recta numérica.", $false, $false, $false, $false,
    $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find the 'Número racional, recta numérica.' paragraph text"
}
$r1.Text = ""

$xml1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00FB4231" w:rsidRPr="000719EE" w:rsidRDefault="00FB4231" w:rsidP="00FB4231"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>n</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">úmero </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>racional,recta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>numérica</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: paragraph "Ubica en cada letra el número correspondiente de
#   acuerdo con su ubicación en la recta numérica." — remove the old
#   "_GoBack" bookmarkStart/bookmarkEnd pair that used to wrap "con ".
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "Ubica en cada letra el número correspondiente de acuerdo con su ubicación en la recta numérica.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find the 'Ubica en cada letra ...' paragraph text"
}
$r2.Text = ""

$xml2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00FB4231" w:rsidRPr="000719EE" w:rsidRDefault="00FB4231" w:rsidP="00FB4231"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">Ubica en cada letra el número correspondiente de acuerdo </w:t></w:r><w:r w:rsidR="00CE6EB1"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">con </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>su ubicación en la recta numérica</w:t></w:r><w:r w:rsidR="000B539E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r2.InsertXML($xml2)

# Restore the original "track changes" document setting (only suppressed it
# above so the two rewrites wouldn't be recorded as w:ins/w:del revisions).
$d.TrackRevisions = $true

Write-Host "Done."
